# Applies the latest cryptos.xlsx price/volume refresh to the active sheet.
# D (Price) / E (Volume(1h)) columns are stored as text in the workbook, and some
# Price values look like plain decimals (e.g. "449.00", "0.320"); a leading
# apostrophe forces Excel to keep them as text (quotePrefix) instead of silently
# coercing them to numbers and losing the formatting (trailing zeros, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.845.77'
$ws.Range("E2").Value = '  -0.99%  '
# Row 3
$ws.Range("D3").Value = '3.407.69'
$ws.Range("E3").Value = '  -0.83%  '
# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.24%  '
# Row 5
$ws.Range("D5").Value = '''409.26'
$ws.Range("E5").Value = '  +0.68%  '
# Row 6
$ws.Range("D6").Value = '''128.68'
$ws.Range("E6").Value = '  -1.39%  '
# Row 7
$ws.Range("E7").Value = '  +6.23%  '
# Row 8
$ws.Range("E8").Value = '  +0.01%  '
# Row 9
$ws.Range("E9").Value = '  +5.37%  '
# Row 10
$ws.Range("D10").Value = '''0.143'
$ws.Range("E10").Value = '  +2.01%  '
# Row 11
$ws.Range("D11").Value = '''43.66'
$ws.Range("E11").Value = '  +2.64%  '
# Row 12
$ws.Range("D12").Value = '''0.0000228'
$ws.Range("E12").Value = '  +37.85%  '
# Row 13
$ws.Range("D13").Value = '''9.33'
$ws.Range("E13").Value = '  +8.19%  '
# Row 14
$ws.Range("E14").Value = '  -0.38%  '
# Row 15
$ws.Range("D15").Value = '''21.36'
$ws.Range("E15").Value = '  +7.14%  '
# Row 16
$ws.Range("D16").Value = '3.946.47'
$ws.Range("E16").Value = '  -0.83%  '
# Row 17
$ws.Range("D17").Value = '3.440.95'
$ws.Range("E17").Value = '  +0.08%  '
# Row 18
$ws.Range("D18").Value = '''12.54'
$ws.Range("E18").Value = '  +8.38%  '
# Row 19
$ws.Range("E19").Value = '  +6.42%  '
# Row 20
$ws.Range("D20").Value = '61.864.86'
$ws.Range("E20").Value = '  -1.09%  '
# Row 21
$ws.Range("D21").Value = '''449.00'
$ws.Range("E21").Value = '  +42.30%  '
# Row 22
$ws.Range("D22").Value = '''91.88'
$ws.Range("E22").Value = '  +8.46%  '
# Row 23
$ws.Range("E23").Value = '  +0.06%  '
# Row 24
$ws.Range("D24").Value = '''13.22'
$ws.Range("E24").Value = '  +2.58%  '
# Row 25
$ws.Range("D25").Value = '''3.30'
$ws.Range("E25").Value = '  +3.62%  '
# Row 26
$ws.Range("D26").Value = '''33.32'
$ws.Range("E26").Value = '  +11.09%  '
# Row 27
$ws.Range("D27").Value = '''9.29'
$ws.Range("E27").Value = '  +13.66%  '
# Row 28
$ws.Range("E28").Value = '  +0.96%  '
# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''2.78'
$ws.Range("E29").Value = '  +2.53%  '
# Row 30
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '''7.63'
$ws.Range("E30").Value = '  -2.12%  '
# Row 31
$ws.Range("D31").Value = '''12.04'
$ws.Range("E31").Value = '  +4.72%  '
# Row 32
$ws.Range("E32").Value = '  -2.36%  '
# Row 33
$ws.Range("E33").Value = '  -0.78%  '
# Row 34
$ws.Range("D34").Value = '''42.54'
$ws.Range("E34").Value = '  -3.85%  '
# Row 36
$ws.Range("E36").Value = '  +3.51%  '
# Row 37
$ws.Range("D37").Value = '''53.97'
$ws.Range("E37").Value = '  +5.09%  '
# Row 38
$ws.Range("D38").Value = '''0.998'
$ws.Range("E38").Value = '  -0.08%  '
# Row 39
$ws.Range("E39").Value = '  +7.77%  '
# Row 40
$ws.Range("E40").Value = '  +1.58%  '
# Row 41
$ws.Range("E41").Value = '  -0.57%  '
# Row 42
$ws.Range("D42").Value = '''0.320'
$ws.Range("E42").Value = '  +0.09%  '
# Row 43
$ws.Range("E43").Value = '  +11.27%  '
# Row 44
$ws.Range("D44").Value = '''143.98'
$ws.Range("E44").Value = '  +0.06%  '
# Row 45
$ws.Range("D45").Value = '''2.60'
$ws.Range("E45").Value = '  +16.44%  '
# Row 46
$ws.Range("E46").Value = '  +0.94%  '
# Row 47
$ws.Range("D47").Value = '''16.70'
$ws.Range("E47").Value = '  -1.72%  '
# Row 48
$ws.Range("D48").Value = '''0.153'
$ws.Range("E48").Value = '  +24.67%  '
# Row 49
$ws.Range("E49").Value = '  +5.53%  '
# Row 50
$ws.Range("D50").Value = '''2.16'
$ws.Range("E50").Value = '  +6.35%  '
# Row 51
$ws.Range("D51").Value = '3.749.95'
$ws.Range("E51").Value = '  -0.74%  '
